$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / mean calculation
$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = -10
$ws.Range("F9").Value = -1
